$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = "shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg , hend_mahmoud@med.asu.edu.eg"
$ws.Range("G9").Value = "wessam.atef@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G15").Value = "nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G27").Value = "nourhanmohamed@med.asu.edu.eg, randa.rabea@med.asu.edu.eg, System, marian.samir@med.asu.edu.eg"
$ws.Range("G29").Value = "neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("G34").Value = "shaimaa.ahmed@med.asu.edu.eg, servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg , hend_mahmoud@med.asu.edu.eg"
$ws.Range("G40").Value = "wessam.atef@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G46").Value = "Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G58").Value = "nourhanmohamed@med.asu.edu.eg, randa.rabea@med.asu.edu.eg, System, marian.samir@med.asu.edu.eg"
$ws.Range("G60").Value = "neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("G65").Value = "Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G71").Value = "Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G76").Value = "nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G77").Value = "Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G91").Value = "neveen.nashaat@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("G95").Value = "servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G96").Value = "servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, nourhan.mahmoud@med.asu.edu.eg "
$ws.Range("G97").Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G108").Value = "nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G120").Value = "nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G122").Value = "neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G126").Value = "servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G127").Value = "servinaz@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, nourhan.mahmoud@med.asu.edu.eg "
$ws.Range("G128").Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
$ws.Range("G138").Value = "nourhanmamdouh@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G139").Value = "Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G140").Value = "Alshimaa_khaled@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg"
$ws.Range("G151").Value = "user@user.com, enas.omran@med.asu.edu.eg, Amr-Saeed@med.asu.edu.eg`n, takyas@med.asu.edu.eg"
$ws.Range("G153").Value = "neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G154").Value = "Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
$ws.Range("G157").Value = "eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G158").Value = "Amira.Sobhy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg"
$ws.Range("G164").Value = "Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"
$ws.Range("G170").Value = "Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"
$ws.Range("G182").Value = "nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G184").Value = "neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G185").Value = "Aya_hamed@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
